$d = $word.ActiveDocument

# The document has a "2SA" Heading2 paragraph immediately followed by a
# paragraph whose only visible content is the italic book title
# "2 Samweli". That whole paragraph (title + its paragraph mark) needs to
# be removed so the "2SA" heading runs straight into the following
# (blank/space) paragraph.
#
# There is a second, unrelated "2 Samweli" paragraph further down (a
# Heading2 for the actual book section) that must NOT be touched, so we
# identify the target paragraph specifically by its italic formatting
# (the only italic run in the whole document).

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "2 Samweli`r" -and $p.Range.Italic -eq -1) {
        $p.Range.Delete()
    }
}
